# Append 8 new NBA game rows (rows 636-643) to Sheet1, mirroring the
# existing table layout:
#   A=Away team  B=Away Pts  C=Home team  D=Home Pts  E=Overtime
#   F=Attend.    G=Arena     H=Win        I=Loss

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Milwaukee Bucks",      122, "Detroit Pistons",        113, "No", 17832, "Little Caesars Arena",      "Milwaukee Bucks",      "Detroit Pistons"),
    @("Cleveland Cavaliers",  126, "Orlando Magic",           99,  "No", 17832, "Amway Center",              "Cleveland Cavaliers",  "Orlando Magic"),
    @("San Antonio Spurs",    123, "Philadelphia 76ers",      133, "No", 17832, "Wells Fargo Center",        "Philadelphia 76ers",   "San Antonio Spurs"),
    @("Memphis Grizzlies",    108, "Toronto Raptors",         100, "No", 17832, "Scotiabank Arena",          "Memphis Grizzlies",    "Toronto Raptors"),
    @("Charlotte Hornets",    128, "Minnesota Timberwolves",  125, "No", 17832, "Target Center",             "Charlotte Hornets",    "Minnesota Timberwolves"),
    @("Boston Celtics",       119, "Dallas Mavericks",        110, "No", 17832, "American Airlines Center",  "Boston Celtics",       "Dallas Mavericks"),
    @("Chicago Bulls",        113, "Phoenix Suns",            115, "No", 17832, "Footprint Center",          "Phoenix Suns",         "Chicago Bulls"),
    @("Atlanta Hawks",        107, "Sacramento Kings",        122, "No", 17832, "Golden 1 Center",           "Sacramento Kings",     "Atlanta Hawks")
)

$startRow = 636
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}

# Move the active selection to match where the user ended up editing.
$ws.Range("G642").Select()
